# Correct the delivered_invitations count for "PSYCHOLOGICAL METHODS".
# A participant who wanted to be removed after prenotification is now
# indicated as a bounce, so delivered_invitations drops from 119 to 118
# (row 17: Source.Title = PSYCHOLOGICAL METHODS).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet 1")
$ws.Range("E17").Value = 118
